# Excel COM-interop edit script
# Summary of changes (per commit message / diff):
#  - Added full commentary (kept header comments, dropped the duplicated
#    "Values" comments on B3/B4 that were accidentally copy-pasted)
#  - Changed the button/motor_control icon from Images/hand.png to
#    Images/button.png (cell B2)
#  - Adjusted visual control parameters for move choice (font color +
#    alignment formatting applied to rows 3 and 4, columns A:C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change the motor_control image reference (B2): hand -> button ---
$ws.Range("B2").Value = "Images/button.png"

# --- Remove the duplicated "Values" comments on B3 and B4 ---
$ws.Range("B3").Comment.Delete()
$ws.Range("B4").Comment.Delete()

# --- Adjust visual formatting for the watch_Tetris / fixation_cross rows ---
$ws.Range("A3:C3").Font.Color = 0
$ws.Range("A4:C4").Font.Color = 0
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B4").HorizontalAlignment = -4131
